# Apply the updated "Current Balance" value for the QANDOVERAL row (row 9)
# Column H holds the current-balance text "<qty>:<status>"; it moves from
# "0:0" to "1:0" in this refreshed report export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H9").Value = "1:0"
